$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 34.111822
$ws.Range("H2").Value = 102.335466
$ws.Range("I2").Value = 0.4228853893909983
$ws.Range("J2").Value = 0.4228853893909983
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 149.829178
$ws.Range("N2").Value = 449.487534
$ws.Range("O2").Value = 0.965236887286734
$ws.Range("P2").Value = 0.965236887286734
$ws.Range("Q2").Value = 5110.946250342316
$ws.Range("R2").Value = 45998.51625308084
$ws.Range("S2").Value = 0.4081845769348056
$ws.Range("T2").Value = 0.4081845769348056

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 34.111822
$ws.Range("H3").Value = 102.335466
$ws.Range("I3").Value = 0.4228853893909983
$ws.Range("J3").Value = 0.4228853893909983
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.021452666666667
$ws.Range("N3").Value = 3.064358
$ws.Range("O3").Value = 0.006580452523633729
$ws.Range("P3").Value = 0.006580452523633729
$ws.Range("Q3").Value = 34.84361154675867
$ws.Range("R3").Value = 313.5925039208281
$ws.Range("S3").Value = 0.002782777227825827
$ws.Range("T3").Value = 0.002782777227825827

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 34.111822
$ws.Range("H4").Value = 102.335466
$ws.Range("I4").Value = 0.4228853893909983
$ws.Range("J4").Value = 0.4228853893909983
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 4.374661666666667
$ws.Range("N4").Value = 13.123985
$ws.Range("O4").Value = 0.02818266018963228
$ws.Range("P4").Value = 0.02818266018963228
$ws.Range("Q4").Value = 149.2276800835567
$ws.Range("R4").Value = 1343.04912075201
$ws.Range("S4").Value = 0.01191803522836684
$ws.Range("T4").Value = 0.01191803522836683

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.88211266666666
$ws.Range("H5").Value = 104.646338
$ws.Range("I5").Value = 0.4324347083490296
$ws.Range("J5").Value = 0.4324347083490295
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 149.829178
$ws.Range("N5").Value = 449.487534
$ws.Range("O5").Value = 0.965236887286734
$ws.Range("P5").Value = 0.965236887286734
$ws.Range("Q5").Value = 5226.358267750054
$ws.Range("R5").Value = 47037.2244097505
$ws.Range("S5").Value = 0.4174019318415639
$ws.Range("T5").Value = 0.4174019318415639

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.88211266666666
$ws.Range("H6").Value = 104.646338
$ws.Range("I6").Value = 0.4324347083490296
$ws.Range("J6").Value = 0.4324347083490295
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.021452666666667
$ws.Range("N6").Value = 3.064358
$ws.Range("O6").Value = 0.006580452523633729
$ws.Range("P6").Value = 0.006580452523633729
$ws.Range("Q6").Value = 35.63042700233378
$ws.Range("R6").Value = 320.673843021004
$ws.Range("S6").Value = 0.002845616067862187
$ws.Range("T6").Value = 0.002845616067862187

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 34.88211266666666
$ws.Range("H7").Value = 104.646338
$ws.Range("I7").Value = 0.4324347083490296
$ws.Range("J7").Value = 0.4324347083490295
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 4.374661666666667
$ws.Range("N7").Value = 13.123985
$ws.Range("O7").Value = 0.02818266018963228
$ws.Range("P7").Value = 0.02818266018963228
$ws.Range("Q7").Value = 152.5974411352145
$ws.Range("R7").Value = 1373.37697021693
$ws.Range("S7").Value = 0.01218716043960344
$ws.Range("T7").Value = 0.01218716043960344

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.67052633333333
$ws.Range("H8").Value = 35.011579
$ws.Range("I8").Value = 0.1446799022599722
$ws.Range("J8").Value = 0.1446799022599721
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 149.829178
$ws.Range("N8").Value = 449.487534
$ws.Range("O8").Value = 0.965236887286734
$ws.Range("P8").Value = 0.965236887286734
$ws.Range("Q8").Value = 1748.585367350687
$ws.Range("R8").Value = 15737.26830615619
$ws.Range("S8").Value = 0.1396503785103644
$ws.Range("T8").Value = 0.1396503785103644

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.67052633333333
$ws.Range("H9").Value = 35.011579
$ws.Range("I9").Value = 0.1446799022599722
$ws.Range("J9").Value = 0.1446799022599721
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.021452666666667
$ws.Range("N9").Value = 3.064358
$ws.Range("O9").Value = 0.006580452523633729
$ws.Range("P9").Value = 0.006580452523633729
$ws.Range("Q9").Value = 11.92089024458689
$ws.Range("R9").Value = 107.288012201282
$ws.Range("S9").Value = 0.000952059227945715
$ws.Range("T9").Value = 0.0009520592279457148

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.67052633333333
$ws.Range("H10").Value = 35.011579
$ws.Range("I10").Value = 0.1446799022599722
$ws.Range("J10").Value = 0.1446799022599721
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 4.374661666666667
$ws.Range("N10").Value = 13.123985
$ws.Range("O10").Value = 0.02818266018963228
$ws.Range("P10").Value = 0.02818266018963228
$ws.Range("Q10").Value = 51.05460418025723
$ws.Range("R10").Value = 459.491437622315
$ws.Range("S10").Value = 0.004077464521662007
$ws.Range("T10").Value = 0.004077464521662006

